$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AMLMono7")

# --- Numeric value edits (filter threshold table) ---

# Row 2
$ws.Range("C2").Value = 10

# Row 3
$ws.Range("L3").Value = 0.9

# Row 4
$ws.Range("L4").Value = 0.7
$ws.Range("N4").Value = 2
$ws.Range("O4").Value = 0.001

# Row 5
$ws.Range("L5").Value = 0.5

# --- Remove the quote-prefix-ish custom font style from the O column notes cells (O2:O5) ---
# These cells previously carried a distinct (black-font) cell style; reset them to the
# workbook's default "Normal" style so no explicit s="2" attribute remains.
$ws.Range("O2").Style = "Normal"
$ws.Range("O3").Style = "Normal"
$ws.Range("O4").Style = "Normal"
$ws.Range("O5").Style = "Normal"

# --- Update the explanatory note text referenced by L8 ---
$ws.Range("L8").Value = "sets maximum fraction of  NVAF relative to TVAF"

# --- Move the active selection ---
$ws.Range("H15").Select() | Out-Null
